# dachser.xlsx update: move the UK (GB / column M) zone assignments that used
# to sit next to the GB postal-code prefixes (column B, rows 2-71) down into
# their own block of rows (102-171) keyed off column A, and clear the old
# column M values from rows 2-71.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Zonen")

# Column B holds the GB postal prefix for rows 2..71, column M the zone
# ("Z0x") assigned to that prefix. Re-home each pair onto a new row
# (100 + original row), with the prefix now in column A and the zone
# still in column M, then clear the original M column cell.
for ($r = 2; $r -le 71; $r++) {
    $prefix = $ws.Cells.Item($r, 2).Value()    # column B = PLZ_GB
    $zone   = $ws.Cells.Item($r, 13).Value()   # column M = GB

    $newRow = $r + 100
    $ws.Cells.Item($newRow, 1).Value = $prefix
    $ws.Cells.Item($newRow, 13).Value = $zone
}

$ws.Range("M2:M71").ClearContents()

# Restore the selection/view state: "Zonen" becomes the active sheet with
# M2:M71 selected (the old GB-zone range, now cleared), scrolled near the
# bottom of the new data block.
$ws.Activate()
$ws.Range("M2:M71").Select()
$excel.ActiveWindow.ScrollRow = 152

# The previously-active sheet ("GWK") no longer keeps the tab-selected flag
# once "Zonen" is activated above, matching the authored file.
